# Fix tests for excel import
#
# The worksheet gains a new "Hyväksymisen ehto" (Acceptance condition)
# column between the existing "Hakemuksentila" (J) and "Vastaanottotila"
# (old K, now L) columns. Inserting a real column lets Excel shift all the
# existing K/L/M data (and styles) one slot to the right automatically, so
# the previous "Vastaanottotila"/"Ilmoittautumistila"/"Julkaistavissa"
# columns land on L/M/N without us having to touch them by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K; existing K/L/M (and their formatting) shift to L/M/N.
$ws.Columns("K").Insert()

# Header for the freshly inserted column; the test data row (row 6) is left
# blank in this column (no acceptance condition for the sample applicant).
$ws.Range("K5").Value = "Hyväksymisen ehto"

# Widen the header columns around the new one so the longer Finnish labels fit.
$ws.Columns("J").ColumnWidth = 13.5
$ws.Columns("K").ColumnWidth = 23.285714285714285
$ws.Columns("L").ColumnWidth = 19.5

# The active selection moves from the old J6 to the new K6.
$ws.Range("K6").Select() | Out-Null
